$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values (as shared-string text) for column A per new row, matching the
# alternating pattern introduced by the diff (rows 251-265).
$colA = @(
    "Login with valid username and password", # 251
    "Login with valid username and password", # 252
    "Create a country",                        # 253
    "Login with valid username and password", # 254
    "Create a country",                        # 255
    "Login with valid username and password", # 256
    "Create a country",                        # 257
    "Login with valid username and password", # 258
    "Create a country",                        # 259
    "Login with valid username and password", # 260
    "Create a country",                        # 261
    "Login with valid username and password", # 262
    "Create a country",                        # 263
    "Login with valid username and password", # 264
    "Create a country"                         # 265
)

$startRow = 251
for ($i = 0; $i -lt $colA.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $colA[$i]
    $ws.Cells.Item($row, 2).Value = "PASSED"
    $ws.Cells.Item($row, 3).Value = "chrome"
}
